# docs: Update QA roadmap with TC-04 test
#
# Adds a 4th test case (TC-04 - "Validar carga de productos") to the
# QA roadmap sheet, following the same row layout used by the other
# three test cases (TC-01..TC-03): one "header" row with the full set
# of test-case metadata (row 15) followed by three "step" continuation
# rows (16-18) that only carry the step/technical-action/expected-result
# columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: clone the number/visual formatting of a reference cell onto a
# target cell via copy/PasteSpecial (xlPasteFormats = -4122). This keeps
# us reusing the workbook's existing style (cellXfs) entries instead of
# minting new ones for every cell.
# ---------------------------------------------------------------------
function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null
}

# ==================== Row 15 (TC-04 header row) ====================
Copy-Format "K14" "A15"
$ws.Range("A15").Value = "TC-04"

Copy-Format "K14" "B15"
$ws.Range("B15").Value = "Validar carga de productos (Happy Path)"

Copy-Format "C2" "C15"
$ws.Range("C15").Value = "P0 (Crítica)"

Copy-Format "D2" "D15"
$ws.Range("D15").Value = "Finalizado"

Copy-Format "E2" "E15"
$ws.Range("E15").Value = '[data-test="login-button"], [data-test="username"], [data-test="password"]'

Copy-Format "F2" "F15"
$ws.Range("F15").Value = "✅ SÍ (Playwright)"

Copy-Format "G2" "G15"
$ws.Range("G15").Value = "PASSED 🟢"

Copy-Format "K14" "H15"
$ws.Range("H15").NumberFormat = "dd/mm/yyyy"
$ws.Range("H15").Value = 46053

Copy-Format "K14" "I15"
$ws.Range("I15").Value = 1

Copy-Format "J2" "J15"
$ws.Range("J15").Value = " Navegar a la página"

Copy-Format "K2" "K15"
$ws.Range("K15").Value = "La página carga correctamente."

Copy-Format "C2" "L15"
$ws.Range("L15").Value = "Hecho"

# ==================== Row 16 (TC-04 step 2) ====================
Copy-Format "K14" "I16"
$ws.Range("I16").Value = 2

Copy-Format "K14" "J16"
$ws.Range("J16").Value = " Iniciar sesion con credenciales validas"

Copy-Format "J2" "K16"
$ws.Range("K16").Value = "Los campos de texto aceptan la entrada."

Copy-Format "C2" "L16"
$ws.Range("L16").Value = "Hecho"

# ==================== Row 17 (TC-04 step 3) ====================
Copy-Format "K14" "I17"
$ws.Range("I17").Value = 3

Copy-Format "K14" "J17"
$ws.Range("J17").Value = " Validar que entramos a pagina Productos"

Copy-Format "K14" "K17"
$ws.Range("K17").Value = "El titulo de la pagina es Products"

Copy-Format "C2" "L17"
$ws.Range("L17").Value = "Hecho"

# ==================== Row 18 (TC-04 step 4) ====================
Copy-Format "K14" "I18"
$ws.Range("I18").Value = 4

Copy-Format "K14" "J18"
$ws.Range("J18").Value = " Validar que se cargan los productos correctamente."

Copy-Format "K14" "K18"
$ws.Range("K18").Value = "Se encuentran 6 productos"

Copy-Format "C2" "L18"
$ws.Range("L18").Value = "Hecho"
